$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save"), matching the header style already used by
# G1 (bold / bordered / centered) so the new column looks consistent with
# the rest of the header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data column H2:H6 filled with 0, left unstyled (same as B2:G6).
$ws.Range("H2:H6").Value = 0
